$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a numeric-looking text value into a cell as a literal string
# (avoids Excel auto-converting "545.43"-style text into a real number).
function Set-TextValue($cell, $val) {
    $cell.Formula = '="' + $val + '"'
    $cell.Copy($null) | Out-Null
    $cell.PasteSpecial(-4163, $null, $null, $null) | Out-Null
}

Set-TextValue $ws.Range('D2') '60.057.18'
$ws.Range('E2').Value = '  -0.33%  '

Set-TextValue $ws.Range('D3') '2.313.70'
$ws.Range('E3').Value = '  -1.94%  '

$ws.Range('E4').Value = '  -0.09%  '

Set-TextValue $ws.Range('D5') '545.43'
$ws.Range('E5').Value = '  -0.21%  '

Set-TextValue $ws.Range('D6') '129.69'
$ws.Range('E6').Value = '  -2.08%  '

$ws.Range('E7').Value = '  -0.03%  '

Set-TextValue $ws.Range('D8') '0.575'
$ws.Range('E8').Value = '  -2.30%  '

Set-TextValue $ws.Range('D9') '2.309.08'
$ws.Range('E9').Value = '  -2.09%  '

$ws.Range('E10').Value = '  -0.12%  '

Set-TextValue $ws.Range('D11') '5.56'
$ws.Range('E11').Value = '  +0.83%  '

$ws.Range('E12').Value = '  -0.59%  '

Set-TextValue $ws.Range('D13') '0.334'
$ws.Range('E13').Value = '  -0.17%  '

Set-TextValue $ws.Range('D14') '23.48'
$ws.Range('E14').Value = '  -1.89%  '

Set-TextValue $ws.Range('D15') '60.001.75'
$ws.Range('E15').Value = '  -0.32%  '

Set-TextValue $ws.Range('D16') '2.720.87'
$ws.Range('E16').Value = '  -2.02%  '

Set-TextValue $ws.Range('D17') '0.0000134'
$ws.Range('E17').Value = '  +0.15%  '

Set-TextValue $ws.Range('D18') '2.306.92'
$ws.Range('E18').Value = '  -1.30%  '

Set-TextValue $ws.Range('D19') '10.54'
$ws.Range('E19').Value = '  -1.71%  '

Set-TextValue $ws.Range('D20') '4.07'
$ws.Range('E20').Value = '  -2.29%  '

Set-TextValue $ws.Range('D21') '313.10'
$ws.Range('E21').Value = '  -0.52%  '

Set-TextValue $ws.Range('D22') '6.55'
$ws.Range('E22').Value = '  -3.89%  '

$ws.Range('E23').Value = '  +0.02%  '

Set-TextValue $ws.Range('D24') '64.11'
$ws.Range('E24').Value = '  +1.27%  '

$ws.Range('E25').Value = '  -0.82%  '

$ws.Range('E26').Value = '  -0.07%  '

Set-TextValue $ws.Range('D27') '7.81'
$ws.Range('E27').Value = '  -1.15%  '

Set-TextValue $ws.Range('D28') '1.37'
$ws.Range('E28').Value = '  +0.70%  '

Set-TextValue $ws.Range('D29') '1.25'
$ws.Range('E29').Value = '  +8.50%  '

Set-TextValue $ws.Range('D30') '171.24'
$ws.Range('E30').Value = '  -0.39%  '

Set-TextValue $ws.Range('D31') '1.72'
$ws.Range('E31').Value = '  -0.53%  '

Set-TextValue $ws.Range('D32') '0.0₃0724'
$ws.Range('E32').Value = '  -0.72%  '

Set-TextValue $ws.Range('D33') '5.97'
$ws.Range('E33').Value = '  +0.31%  '

$ws.Range('B34').Value = 'ImmutableX'
$ws.Range('C34').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D34') '1.36'
$ws.Range('E34').Value = '  -4.14%  '

$ws.Range('B35').Value = 'PolygonEcosystemToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue $ws.Range('D35') '0.381'
$ws.Range('E35').Value = '  -0.64%  '

$ws.Range('B36').Value = 'EthereumClassic'
$ws.Range('C36').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue $ws.Range('D36') '17.92'
$ws.Range('E36').Value = '  -0.82%  '

$ws.Range('B37').Value = 'USDe'
$ws.Range('C37').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
Set-TextValue $ws.Range('D37') '0.999'
$ws.Range('E37').Value = '  -0.01%  '

Set-TextValue $ws.Range('D38') '0.999'
$ws.Range('E38').Value = '  -0.09%  '

Set-TextValue $ws.Range('D39') '4.06'
$ws.Range('E39').Value = '  -2.53%  '

Set-TextValue $ws.Range('D40') '315.41'
$ws.Range('E40').Value = '  -2.21%  '

Set-TextValue $ws.Range('D41') '37.97'
$ws.Range('E41').Value = '  -0.41%  '

$ws.Range('E42').Value = '  -1.46%  '

Set-TextValue $ws.Range('D43') '137.39'
$ws.Range('E43').Value = '  -3.73%  '

Set-TextValue $ws.Range('D44') '3.49'
$ws.Range('E44').Value = '  +0.64%  '

Set-TextValue $ws.Range('D45') '0.0941'
$ws.Range('E45').Value = '  -1.01%  '

Set-TextValue $ws.Range('D46') '19.03'
$ws.Range('E46').Value = '  -2.08%  '

Set-TextValue $ws.Range('D47') '0.564'
$ws.Range('E47').Value = '  -0.23%  '

$ws.Range('E48').Value = '  -1.31%  '

Set-TextValue $ws.Range('D49') '0.0215'
$ws.Range('E49').Value = '  -0.08%  '

Set-TextValue $ws.Range('D50') '0.0₆0217'
$ws.Range('E50').Value = '  +2.32%  '

$ws.Range('B51').Value = 'WhiteBITCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
Set-TextValue $ws.Range('D51') '10.91'
$ws.Range('E51').Value = '  -0.88%  '

$excel.CutCopyMode = $false